# Append the 01/08/2026 profit-data row to the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet currently has data through row 44 (A1:J44) -> new row goes at 45.
$row = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

# Column A holds text dates (e.g. "01/07/2026"), not real Excel dates.
# Force text formatting BEFORE assigning the value so Excel doesn't
# auto-convert the slash-separated string into a date serial number,
# then drop the formatting override again so the cell ends up with no
# explicit style, matching the rest of the column.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01/08/2026"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value  = 12635.58
$ws.Cells.Item($row, 3).Value  = 0.2125815996498773
$ws.Cells.Item($row, 4).Value  = 0.7874184003501227
$ws.Cells.Item($row, 5).Value  = -126.25
$ws.Cells.Item($row, 6).Value  = -21.58
$ws.Cells.Item($row, 7).Value  = -20692.27
$ws.Cells.Item($row, 8).Value  = -67.53
$ws.Cells.Item($row, 9).Value  = -416.76
$ws.Cells.Item($row, 10).Value = -13.43
